$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the data (all columns except the leading id column A) between
#        row 2 and row 3 ---
$rng2 = $ws.Range("B2:AC2")
$rng3 = $ws.Range("B3:AC3")
$vals2 = $rng2.Value()
$vals3 = $rng3.Value()
$rng2.Value = $vals3
$rng3.Value = $vals2

# --- 2) Swap the data (all columns except the leading id column A) between
#        row 187 and row 188 ---
$rng187 = $ws.Range("B187:AC187")
$rng188 = $ws.Range("B188:AC188")
$vals187 = $rng187.Value()
$vals188 = $rng188.Value()
$rng187.Value = $vals188
$rng188.Value = $vals187

# --- 3) Refresh the two still-upcoming fixtures: row 267 picks up what used
#        to be row 273's data, row 268 picks up what used to be row 274's
#        data (id in column A is left untouched) ---
$ws.Range("B267").Value = 7979356
$ws.Range("E267").Value = 45403.5625
$ws.Range("F267").Value = "Union Saint Gilloise"
$ws.Range("G267").Value = "Club Brugge"
$ws.Range("K267").Value = 2.3
$ws.Range("L267").Value = 3.4
$ws.Range("M267").Value = 3
$ws.Range("N267").Value = 2.3
$ws.Range("O267").Value = 3.4
$ws.Range("P267").Value = 3
$ws.Range("Q267").Value = -0.25
$ws.Range("R267").Value = 2
$ws.Range("S267").Value = 1.85
$ws.Range("T267").Value = 2.5
$ws.Range("U267").Value = 1.85
$ws.Range("V267").Value = 2

$ws.Range("B268").Value = 7979476
$ws.Range("E268").Value = 45403.59375
$ws.Range("F268").Value = "Charleroi"
$ws.Range("G268").Value = "Eupen"
$ws.Range("K268").Value = 1.666
$ws.Range("L268").Value = 3.75
$ws.Range("M268").Value = 5
$ws.Range("N268").Value = 1.666
$ws.Range("O268").Value = 3.75
$ws.Range("P268").Value = 5
$ws.Range("Q268").Value = -0.75
$ws.Range("R268").Value = 1.875
$ws.Range("S268").Value = 1.975
$ws.Range("T268").Value = 2.5
$ws.Range("U268").Value = 1.875
$ws.Range("V268").Value = 1.975

# --- 4) Drop the six fixtures that are no longer part of the feed (old rows
#        269-274); the sheet's used range shrinks from AC274 to AC268 ---
$ws.Rows("269:274").Delete()
